$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume snapshot with the latest
# scrape. Cells D/E are stored as plain text in the source sheet (not
# numbers), so numeric-looking values (e.g. "618.49", "1.00", "0.0000131")
# are entered with a leading apostrophe to force Excel to keep them as text
# instead of silently re-parsing them as numbers (which would drop
# significant trailing/leading zeros and change the stored precision).
# The Style reset afterwards keeps the cell's default (unstyled) appearance
# instead of picking up an incidental "Text" number format.

$ws.Range("D2").Value = '68.184.66'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '3.562.62'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'618.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").Value = "'154.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").Value = '3.562.04'
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +2.21%  '
$ws.Range("E10").Value = '  +5.23%  '
$ws.Range("D11").Value = "'7.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.90%  '
$ws.Range("E12").Value = '  +4.05%  '
$ws.Range("D13").Value = "'33.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.35%  '
$ws.Range("E14").Value = '  +1.33%  '
$ws.Range("D15").Value = '4.168.63'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").Value = '3.565.20'
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").Value = '68.281.92'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").Value = "'0.116"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("E19").Value = '  +5.24%  '
$ws.Range("D20").Value = "'15.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.53%  '
$ws.Range("D21").Value = "'10.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.68%  '
$ws.Range("D22").Value = "'454.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("D23").Value = "'0.644"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.02%  '
$ws.Range("D24").Value = "'78.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = "'0.0000131"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("D26").Value = '3.708.60'
$ws.Range("E26").Value = '  +1.78%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = "'9.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.64%  '
$ws.Range("D29").Value = "'10.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.95%  '
$ws.Range("D30").Value = "'1.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.45%  '
$ws.Range("E31").Value = '  +3.45%  '
$ws.Range("D32").Value = "'0.170"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.30%  '
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +4.87%  '
$ws.Range("D35").Value = "'26.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").Value = '  +4.79%  '
$ws.Range("D37").Value = '3.559.20'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").Value = "'8.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.28%  '
$ws.Range("D39").Value = "'2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.80%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = "'181.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.95%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = "'0.0917"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.84%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = "'5.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.17%  '
$ws.Range("D45").Value = "'30.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.22%  '
$ws.Range("D46").Value = "'0.897"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("D47").Value = "'46.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.11%  '
$ws.Range("E48").Value = '  +4.85%  '
$ws.Range("D49").Value = "'2.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.69%  '
$ws.Range("D50").Value = "'7.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.51%  '
$ws.Range("D51").Value = "'0.262"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.93%  '
